$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.221.11"

$ws.Range("D3").Value = "1.857.78"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6984"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07725"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.82%  "

$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.26"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08164"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("D12").Value = "1.823.64"
$ws.Range("E12").Value = "  -3.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7168"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.153"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").Value = "29.222.13"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.752"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007733"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.77"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").Value = "2.111.34"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9994"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.430"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1483"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.62"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.005"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.044"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.416"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.431"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.483"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.021"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05188"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.166"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7078"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.656"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01845"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.722"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9384"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.56%  "

$ws.Range("D42").Value = "1.138.98"
$ws.Range("E42").Value = "  +8.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4274"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.33"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.791"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.11%  "

$ws.Range("D49").Value = "2.008.52"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.139"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.948"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.93%  "

# Row 44 and 45: Aave/FraxShare swap with updated data
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.883"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.68"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.38%  "
